$d = $word.ActiveDocument

# The document currently holds several one-line "transition" paragraphs
# (Aaa / bbb / a#bb / bb#aa / <empty> / gabriel / ab#). The edit collapses
# all of that down to a single paragraph describing the adjusted
# transition, reusing the first paragraph's formatting.

# 1. Delete every paragraph after the first one.
if ($d.Paragraphs.Count -gt 1) {
    $deleteRange = $d.Range($d.Paragraphs.Item(2).Range.Start, `
                             $d.Paragraphs.Item($d.Paragraphs.Count).Range.End)
    $deleteRange.Delete()
}

# 2. Replace the remaining paragraph's text with the new value.
$d.Content.Find.Execute("Aaa", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "b#b__", 2)
